$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 40 data rows (rows 2-41, below the header in row 1).
# The update adds one more weekly observation (new "Primera"/"Segunda" pair) at
# the top of the data and keeps the rest of the history, so every existing
# data row shifts down by two rows (2->4 .. 41->43).

# Make sure the two rows that will land beyond the old used range (42:43)
# already carry the same date format as the rest of column D before we paste
# into them - this avoids Excel minting a stray "short date" style for cells
# outside the previous used range.
$ws.Range("D42:D43").NumberFormat = $ws.Range("D41").NumberFormat

# Shift the existing data rows 2:41 down to 4:43 (values + formatting).
$ws.Range("A2:R41").Copy() | Out-Null
$ws.Range("A4:R43").PasteSpecial(-4104) | Out-Null
$excel.CutCopyMode = 0

# Rows 2 and 3 keep their original formatting (they were not touched by the
# paste above); overwrite just the cells whose values actually change for
# this new week's observation.
$ws.Range("D2").Value2 = 45160
$ws.Range("J2").Value2 = 54
$ws.Range("K2").Value2 = 23000
$ws.Range("L2").Value2 = 23000
$ws.Range("M2").Value2 = 23000
$ws.Range("P2").Value2 = 1533

$ws.Range("D3").Value2 = 45160
$ws.Range("J3").Value2 = 50
$ws.Range("K3").Value2 = 17000
$ws.Range("L3").Value2 = 17000
$ws.Range("M3").Value2 = 17000
$ws.Range("P3").Value2 = 1133

Write-Output "done"
